$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 183, shifting the existing rows 183:199 down to 184:200
$ws.Rows.Item(183).Insert()

# Populate the newly inserted row 183 with the new weekly record
$ws.Range("A183").Value = 5
$ws.Range("B183").Value = "Macroferia Regional de Talca"
$ws.Range("C183").Value = "Maule"
$ws.Range("D183").Value = 44578
$ws.Range("E183").Value = 7
$ws.Range("F183").Value = "Fruta"
$ws.Range("G183").Value = 100108
$ws.Range("H183").Value = "Tropicales y subtropicales"
$ws.Range("I183").Value = 100108005
$ws.Range("J183").Value = "Piña"
$ws.Range("K183").Value = "Caramelo"
$ws.Range("L183").Value = "Segunda"
$ws.Range("M183").Value = 250
$ws.Range("N183").Value = 14000
$ws.Range("O183").Value = 14000
$ws.Range("P183").Value = 14000
$ws.Range("Q183").Value = "`$/caja 14 unidades"
$ws.Range("R183").Value = "Ecuador"
$ws.Range("S183").Value = 1000
$ws.Range("T183").Value = 14
